# Auto-applies the diff to both '展览' (sheet1) and '全部类型' (sheet4) worksheets,
# which are identical in content before this edit and receive identical edits.
$wb = $excel.ActiveWorkbook

$sheetNames = @('展览', '全部类型')

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Simple 'want-to-go count' (F column) bumps on unchanged rows ---
    $ws.Cells.Item(2, 6).Value = 1555
    $ws.Cells.Item(3, 6).Value = 46
    $ws.Cells.Item(4, 6).Value = 1017
    $ws.Cells.Item(5, 6).Value = 18
    $ws.Cells.Item(7, 6).Value = 2587
    $ws.Cells.Item(9, 6).Value = 1622
    $ws.Cells.Item(15, 6).Value = 51

    # --- Rows 10-14 replaced with updated event listings; rows 16-17 appended new ---

    # Row 10
    $ws.Cells.Item(10, 1).Value = 9
    $ws.Cells.Item(10, 2).Value = '''2024.03.23'
    $ws.Cells.Item(10, 3).Value = '南昌·运动番only春季集训（取消）'
    $ws.Cells.Item(10, 4).Value = '创新三路777号 南昌小飞侠章鱼文化体育公园'
    $ws.Cells.Item(10, 5).Value = '''2024.03.23 10:00-03.24 17:00'
    $ws.Cells.Item(10, 6).Value = 182
    $ws.Cells.Item(10, 7).Value = '不可售'
    $ws.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81950'
    $ws.Cells.Item(10, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/bm4uH4qB1708425538357.jpeg'

    # Row 11
    $ws.Cells.Item(11, 1).Value = 10
    $ws.Cells.Item(11, 2).Value = '''2024.03.24'
    $ws.Cells.Item(11, 3).Value = '南昌·AP动漫游戏  嘉年华内场票-小N&子音'
    $ws.Cells.Item(11, 4).Value = '八一桥街道青山南路118号 蓝海会展中心'
    $ws.Cells.Item(11, 5).Value = '''2024.03.24 09:00-03.24 17:00'
    $ws.Cells.Item(11, 6).Value = 65
    $ws.Cells.Item(11, 7).Value = 218
    $ws.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81973'
    $ws.Cells.Item(11, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/zbG5HICL1708504962467.jpeg'

    # Row 12
    $ws.Cells.Item(12, 1).Value = 11
    $ws.Cells.Item(12, 2).Value = '''2024.03.30'
    $ws.Cells.Item(12, 3).Value = '南昌·CM01动漫游戏博览会'
    $ws.Cells.Item(12, 4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
    $ws.Cells.Item(12, 5).Value = '''2024.03.30 10:00-03.31 17:00'
    $ws.Cells.Item(12, 6).Value = 525
    $ws.Cells.Item(12, 7).Value = 55
    $ws.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81691'
    $ws.Cells.Item(12, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/9cMJMElF1708938074308.png'

    # Row 13
    $ws.Cells.Item(13, 1).Value = 12
    $ws.Cells.Item(13, 2).Value = '''2024.03.30'
    $ws.Cells.Item(13, 3).Value = '鹰潭·原×铁×崩only'
    $ws.Cells.Item(13, 4).Value = '南站路24号 回禾酒店(鹰潭火车站店)'
    $ws.Cells.Item(13, 5).Value = '''2024.03.30 10:00-03.30 17:00'
    $ws.Cells.Item(13, 6).Value = 27
    $ws.Cells.Item(13, 7).Value = 60
    $ws.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81097'
    $ws.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg'

    # Row 14
    $ws.Cells.Item(14, 1).Value = 13
    $ws.Cells.Item(14, 2).Value = '''2024.03.31'
    $ws.Cells.Item(14, 3).Value = '新余·文旅国漫嘉年华暨BM次元盛典'
    $ws.Cells.Item(14, 4).Value = '五一南路与仙女湖大道交叉口西北 老上海风情街白金汉宫'
    $ws.Cells.Item(14, 5).Value = '''2024.03.31 10:00-03.31 17:00'
    $ws.Cells.Item(14, 6).Value = 3
    $ws.Cells.Item(14, 7).Value = 60
    $ws.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82208'
    $ws.Cells.Item(14, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/UOMzhQHg1709202735831.png'

    # Row 16
    $ws.Cells.Item(16, 1).Value = 15
    $ws.Cells.Item(16, 2).Value = '''2024.04.13'
    $ws.Cells.Item(16, 3).Value = '南昌·原X穹X崩only'
    $ws.Cells.Item(16, 4).Value = '丰和北大道299号 新吉花园酒店'
    $ws.Cells.Item(16, 5).Value = '''2024.04.13 10:00-04.13 17:00'
    $ws.Cells.Item(16, 6).Value = 76
    $ws.Cells.Item(16, 7).Value = 65
    $ws.Cells.Item(16, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80807'
    $ws.Cells.Item(16, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/kfK13XvH1709202705153.jpeg'

    # Row 17
    $ws.Cells.Item(17, 1).Value = 16
    $ws.Cells.Item(17, 2).Value = '''2024.04.13'
    $ws.Cells.Item(17, 3).Value = '南昌·第二届漫拥动漫嘉年华mini'
    $ws.Cells.Item(17, 4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
    $ws.Cells.Item(17, 5).Value = '''2024.04.13 10:00-04.14 18:00'
    $ws.Cells.Item(17, 6).Value = 6
    $ws.Cells.Item(17, 7).Value = 39.9
    $ws.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82210'
    $ws.Cells.Item(17, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/KYd0bfk11709203777701.png'

    # New rows 16 and 17 need the bold/bordered/centered style in column A,
    # matching the existing style used by A2:A15. Copy format from A15 (unaffected by the
    # value-only writes above) onto A16/A17 without touching the values just written.
    $ws.Range('A15').Copy()
    $ws.Range('A16').PasteSpecial(-4122)
    $ws.Range('A17').PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

